$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.022634547218986
$ws.Range("D2").Value = 1.027289883076897
$ws.Range("E2").Value = 1.032823614672814
$ws.Range("F2").Value = 1.042712067092526
$ws.Range("I2").Value = 1.028440351683639
$ws.Range("J2").Value = 1.027819279863118
$ws.Range("K2").Value = 1.030110043985621
$ws.Range("L2").Value = 1.035627733886319
$ws.Range("M2").Value = 1.045487981023039
$ws.Range("N2").Value = 1.013320039186339

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.023524801313748
$ws.Range("D3").Value = 1.027929781129339
$ws.Range("E3").Value = 1.033688454792816
$ws.Range("F3").Value = 1.043829514048163
$ws.Range("I3").Value = 1.028540745559051
$ws.Range("J3").Value = 1.028347975078077
$ws.Range("K3").Value = 1.030558089848618
$ws.Range("L3").Value = 1.036301276812463
$ws.Range("M3").Value = 1.046415512701386
$ws.Range("N3").Value = 1.013496837698308

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.024101291327291
$ws.Range("D4").Value = 1.028343903086897
$ws.Range("E4").Value = 1.034248896023099
$ws.Range("F4").Value = 1.044553918907022
$ws.Range("I4").Value = 1.028604341663883
$ws.Range("J4").Value = 1.028689896723169
$ws.Range("K4").Value = 1.030847384713133
$ws.Range("L4").Value = 1.036737297164336
$ws.Range("M4").Value = 1.047016442363344
$ws.Range("N4").Value = 1.01361112520472

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.024343751004691
$ws.Range("D5").Value = 1.02851801399006
$ws.Range("E5").Value = 1.034484703398088
$ws.Range("F5").Value = 1.044858778967335
$ws.Range("I5").Value = 1.028630749995278
$ws.Range("J5").Value = 1.028833596621795
$ws.Range("K5").Value = 1.030968854536975
$ws.Range("L5").Value = 1.036920645320144
$ws.Range("M5").Value = 1.047269252853689
$ws.Range("N5").Value = 1.013659144250849

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.024384467074947
$ws.Range("D6").Value = 1.028547248754279
$ws.Range("E6").Value = 1.034524308061161
$ws.Range("F6").Value = 1.044909985063761
$ws.Range("I6").Value = 1.028635164850385
$ws.Range("J6").Value = 1.028857721860557
$ws.Range("K6").Value = 1.030989241041538
$ws.Range("L6").Value = 1.036951432903376
$ws.Range("M6").Value = 1.047311711353766
$ws.Range("N6").Value = 1.013667205240912

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.024104530680965
$ws.Range("D7").Value = 1.028346229512136
$ws.Range("E7").Value = 1.034252046116826
$ws.Range("F7").Value = 1.044557991204929
$ws.Range("I7").Value = 1.028604695820764
$ws.Range("J7").Value = 1.02869181702222
$ws.Range("K7").Value = 1.030849008388643
$ws.Range("L7").Value = 1.036739746895032
$ws.Range("M7").Value = 1.047019819724797
$ws.Range("N7").Value = 1.01361176694579

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.022935321787505
$ws.Range("D8").Value = 1.027506125197307
$ws.Range("E8").Value = 1.033115718442022
$ws.Range("F8").Value = 1.043089435052936
$ws.Range("I8").Value = 1.028474562639299
$ws.Range("J8").Value = 1.027997991462755
$ws.Range("K8").Value = 1.030261591041546
$ws.Range("L8").Value = 1.035855320290934
$ws.Range("M8").Value = 1.045801287837854
$ws.Range("N8").Value = 1.01337981223024

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.02087841586815
$ws.Range("D9").Value = 1.026026321666619
$ws.Range("E9").Value = 1.031119790602394
$ws.Range("F9").Value = 1.040511973912917
$ws.Range("I9").Value = 1.028234819296704
$ws.Range("J9").Value = 1.026774050828654
$ws.Range("K9").Value = 1.02922178185718
$ws.Range("L9").Value = 1.034298370870074
$ws.Range("M9").Value = 1.043659902962463
$ws.Range("N9").Value = 1.012970230078871

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.019509490108021
$ws.Range("D10").Value = 1.025040260007554
$ws.Range("E10").Value = 1.029793567471242
$ws.Range("F10").Value = 1.038800665593501
$ws.Range("I10").Value = 1.028068010858298
$ws.Range("J10").Value = 1.02595725338935
$ws.Range("K10").Value = 1.028525478954284
$ws.Range("L10").Value = 1.033261488449885
$ws.Range("M10").Value = 1.042236287473016
$ws.Range("N10").Value = 1.012696627407784

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.018917299873852
$ws.Range("D11").Value = 1.024613416231746
$ws.Range("E11").Value = 1.029220356217333
$ws.Range("F11").Value = 1.038061322283728
$ws.Range("I11").Value = 1.027994132682564
$ws.Range("J11").Value = 1.025603384627842
$ws.Range("K11").Value = 1.028223251467916
$ws.Range("L11").Value = 1.032812776354947
$ws.Range("M11").Value = 1.041620800957874
$ws.Range("N11").Value = 1.01257802915268

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.018697419786595
$ws.Range("D12").Value = 1.024454888141358
$ws.Range("E12").Value = 1.029007599324792
$ws.Range("F12").Value = 1.037786948248804
$ws.Range("I12").Value = 1.027966443844959
$ws.Range("J12").Value = 1.025471914403315
$ws.Range("K12").Value = 1.028110883082952
$ws.Range("L12").Value = 1.032646145458578
$ws.Range("M12").Value = 1.041392325272956
$ws.Range("N12").Value = 1.012533957844252

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.018744580915209
$ws.Range("D13").Value = 1.024488891993645
$ws.Range("E13").Value = 1.029053229172656
$ws.Range("F13").Value = 1.037845791038045
$ws.Range("I13").Value = 1.027972394371771
$ws.Range("J13").Value = 1.025500116446986
$ws.Range("K13").Value = 1.028134991336359
$ws.Range("L13").Value = 1.032681886491727
$ws.Range("M13").Value = 1.041441327578349
$ws.Range("N13").Value = 1.012543412126917

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.018899122762231
$ws.Range("D14").Value = 1.024600311819264
$ws.Range("E14").Value = 1.029202766409702
$ws.Range("F14").Value = 1.038038637314959
$ws.Range("I14").Value = 1.027991848954087
$ws.Range("J14").Value = 1.025592517815837
$ws.Range("K14").Value = 1.028213965251663
$ws.Range("L14").Value = 1.032799001756097
$ws.Range("M14").Value = 1.041601912155558
$ws.Range("N14").Value = 1.012574386582615

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.018994352486651
$ws.Range("D15").Value = 1.024668964039451
$ws.Range("E15").Value = 1.02929492238075
$ws.Range("F15").Value = 1.038157489557066
$ws.Range("I15").Value = 1.028003802825318
$ws.Range("J15").Value = 1.025649445718494
$ws.Range("K15").Value = 1.028262609456707
$ws.Range("L15").Value = 1.032871165765082
$ws.Range("M15").Value = 1.041700872644748
$ws.Range("N15").Value = 1.012593468508632

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.019548803788632
$ws.Range("D16").Value = 1.025068591040205
$ws.Range("E16").Value = 1.029831631885361
$ws.Range("F16").Value = 1.038849768521136
$ws.Range("I16").Value = 1.028072879213731
$ws.Range("J16").Value = 1.025980734569653
$ws.Range("K16").Value = 1.028545521642261
$ws.Range("L16").Value = 1.03329127365602
$ws.Range("M16").Value = 1.042277155346394
$ws.Range("N16").Value = 1.012704495749757

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.019896748088348
$ws.Range("D17").Value = 1.025319301710286
$ws.Range("E17").Value = 1.030168578120179
$ws.Range("F17").Value = 1.039284463102742
$ws.Range("I17").Value = 1.02811576790777
$ws.Range("J17").Value = 1.026188493091866
$ws.Range("K17").Value = 1.02872279198949
$ws.Range("L17").Value = 1.033554867822222
$ws.Range("M17").Value = 1.042638896822823
$ws.Range("N17").Value = 1.012774106599085

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.020099752425355
$ws.Range("D18").Value = 1.025465549193886
$ws.Range("E18").Value = 1.030365214589977
$ws.Range("F18").Value = 1.039538173634442
$ws.Range("I18").Value = 1.028140625013646
$ws.Range("J18").Value = 1.026309656719346
$ws.Range("K18").Value = 1.028826120880698
$ws.Range("L18").Value = 1.033708643327572
$ws.Range("M18").Value = 1.042849985948914
$ws.Range("N18").Value = 1.012814697200349

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.020168980838826
$ws.Range("D19").Value = 1.025515417827153
$ws.Range("E19").Value = 1.030432279712091
$ws.Range("F19").Value = 1.039624709571219
$ws.Range("I19").Value = 1.028149073634025
$ws.Range("J19").Value = 1.026350967234815
$ws.Range("K19").Value = 1.028861341472292
$ws.Range("L19").Value = 1.033761081094811
$ws.Range("M19").Value = 1.042921977368085
$ws.Range("N19").Value = 1.012828535460896

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.019859411332359
$ws.Range("D20").Value = 1.025292401555875
$ws.Range("E20").Value = 1.03013241647761
$ws.Range("F20").Value = 1.039237807876296
$ws.Range("I20").Value = 1.028111182813779
$ws.Range("J20").Value = 1.026166204465346
$ws.Range("K20").Value = 1.028703779781404
$ws.Range("L20").Value = 1.033526584019489
$ws.Range("M20").Value = 1.042600075920772
$ws.Range("N20").Value = 1.012766639276024

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.018853611665032
$ws.Range("D21").Value = 1.024567500874034
$ws.Range("E21").Value = 1.029158727013525
$ws.Range("F21").Value = 1.037981841957018
$ws.Range("I21").Value = 1.027986126883801
$ws.Range("J21").Value = 1.025565308670331
$ws.Range("K21").Value = 1.028190712364592
$ws.Range("L21").Value = 1.032764513096966
$ws.Range("M21").Value = 1.041554620033711
$ws.Range("N21").Value = 1.012565265888915

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.01822172204166
$ws.Range("D22").Value = 1.024111847532657
$ws.Range("E22").Value = 1.028547452007058
$ws.Range("F22").Value = 1.03719361894257
$ws.Range("I22").Value = 1.027906069344612
$ws.Range("J22").Value = 1.025187342357953
$ws.Range("K22").Value = 1.02786750428558
$ws.Range("L22").Value = 1.032285605134816
$ws.Range("M22").Value = 1.040898130619859
$ws.Range("N22").Value = 1.012438546808549

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.018556651267987
$ws.Range("D23").Value = 1.02435338605305
$ws.Range("E23").Value = 1.028871412631021
$ws.Range("F23").Value = 1.037611332877063
$ws.Range("I23").Value = 1.027948644715302
$ws.Range("J23").Value = 1.025387724232727
$ws.Range("K23").Value = 1.028038901619062
$ws.Range("L23").Value = 1.032539460587664
$ws.Range("M23").Value = 1.041246069075467
$ws.Range("N23").Value = 1.012505733070722

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.019876282037506
$ws.Range("D24").Value = 1.025304556540251
$ws.Range("E24").Value = 1.030148756054033
$ws.Range("F24").Value = 1.039258888868886
$ws.Range("I24").Value = 1.028113255112173
$ws.Range("J24").Value = 1.026176275793216
$ws.Range("K24").Value = 1.028712370795964
$ws.Range("L24").Value = 1.033539364173808
$ws.Range("M24").Value = 1.042617617132891
$ws.Range("N24").Value = 1.012770013475711

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.02140976689525
$ws.Range("D25").Value = 1.026408809553403
$ws.Range("E25").Value = 1.031635016525695
$ws.Range("F25").Value = 1.04117708009192
$ws.Range("I25").Value = 1.028298031407603
$ws.Range("J25").Value = 1.02709062060638
$ws.Range("K25").Value = 1.02949114814499
$ws.Range("L25").Value = 1.034700692956171
$ws.Range("M25").Value = 1.044212806076294
$ws.Range("N25").Value = 1.013076215097097

